# Actualizacion automatica 2025-07-24 10:00:08
# Updates the monthly sales figures after a new sale was recorded for
# ALMEIDA CUATIN JHONATHANN CARLOS's clients:
#   - AUCANSHALA ALLAICA FREDDY HERNAN: +142.56 in "PIEDRA SINTERIZADA"
#   - HERRERA CAICEDO LUIS FRANKLIN:    +63.06  in "PORCELANATO"
# which ripples through the per-group sheet, the monthly sheet and the
# monthly-compliance summary sheet.

$wb = $excel.ActiveWorkbook

$wsGrupo       = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual     = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- VENTAS POR GRUPO ---------------------------------------------------
$wsGrupo.Range("L4").Value = 142.56
$wsGrupo.Range("M14").Value = 374.03
$wsGrupo.Range("L32").Value = "1 de 30"

# --- VENTA MENSUAL --------------------------------------------------------
$wsMensual.Range("F4").Value = 1190.78
$wsMensual.Range("F14").Value = 374.03
$wsMensual.Range("F32").Value = 8793.26

# --- CUMPLIMIENTO MENSUAL --------------------------------------------------
$wsCumplimiento.Range("D14").Value = 142.56
$wsCumplimiento.Range("E14").Value = 384.47
$wsCumplimiento.Range("F14").Value = 0.2704969356583117

$wsCumplimiento.Range("D15").Value = 5766.91
$wsCumplimiento.Range("E15").Value = 17691.91
$wsCumplimiento.Range("F15").Value = 0.2458312054911543

$wsCumplimiento.Range("D18").Value = 8782.98
$wsCumplimiento.Range("E18").Value = 25151.73607548726
$wsCumplimiento.Range("F18").Value = 0.2588199052693529
